$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# D3: "iaest-measure:tipo-de-hogar" -> "iaest-dimension:tipo-de-hogar"
$ws.Range("D3").Value = "iaest-dimension:tipo-de-hogar"

# D4: "medida" -> "dim"
$ws.Range("D4").Value = "dim"

# D5: "xsd:string" -> "skos:Concept"
$ws.Range("D5").Value = "skos:Concept"

# New row 6, cell D6: "mapping-tipo-de-hogar.xlsx" (copy D5's formatting first, then set the value)
$ws.Range("D5").Copy()
$ws.Range("D6").PasteSpecial(-4122)
$ws.Range("D6").Value = "mapping-tipo-de-hogar.xlsx"
